$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): B1:L1 - mean/std pairs for each horizon -----------
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "One Year Base mean"
$ws.Range("D1").Value = "One Year Base std"
$ws.Range("E1").Value = "Two Year Base mean"
$ws.Range("F1").Value = "Two Year Base std"
$ws.Range("G1").Value = "Three Year Base mean"
$ws.Range("H1").Value = "Three Year Base std"
$ws.Range("I1").Value = "Five Year Base mean"
$ws.Range("J1").Value = "Five Year Base std"
$ws.Range("K1").Value = "Ten Year Base mean"
$ws.Range("L1").Value = "Ten Year Base std"

# New header cells (H1:L1) need the same style as the existing header cells
$ws.Range("G1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: LR ---------------------------------------------------------
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.9041850056937051
$ws.Range("D2").Value = 0.01114962016910395
$ws.Range("E2").Value = 0.9015344564248444
$ws.Range("F2").Value = 0.01113865571179588
$ws.Range("G2").Value = 0.8993876526275825
$ws.Range("H2").Value = 0.01641469494757438
$ws.Range("I2").Value = 0.8970284126711349
$ws.Range("J2").Value = 0.01472074846364865
$ws.Range("K2").Value = 0.8915626060400408
$ws.Range("L2").Value = 0.009703811068141943

# --- Row 3: LDA ----------------------------------------------------------
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.9066248245544344
$ws.Range("D3").Value = 0.009890860256044244
$ws.Range("E3").Value = 0.8981620590389563
$ws.Range("F3").Value = 0.01330489627007841
$ws.Range("G3").Value = 0.892390971550341
$ws.Range("H3").Value = 0.01255113168432723
$ws.Range("I3").Value = 0.8862411926447843
$ws.Range("J3").Value = 0.01945527763094153
$ws.Range("K3").Value = 0.8891873091279268
$ws.Range("L3").Value = 0.0109814743274516

# --- Row 4: KNN ------------------------------------------------------------
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.8672551044728689
$ws.Range("D4").Value = 0.006661247753989168
$ws.Range("E4").Value = 0.8718602551654827
$ws.Range("F4").Value = 0.01450145360826683
$ws.Range("G4").Value = 0.8720867573757241
$ws.Range("H4").Value = 0.01671627558821652
$ws.Range("I4").Value = 0.8777231196654638
$ws.Range("J4").Value = 0.01467693970183832
$ws.Range("K4").Value = 0.8787184707612262
$ws.Range("L4").Value = 0.009236432685052816

# --- Row 5: DTREE (renamed from CART) --------------------------------------
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.8932827520457615
$ws.Range("D5").Value = 0.01121723230674803
$ws.Range("E5").Value = 0.8818082454676046
$ws.Range("F5").Value = 0.008798122488490578
$ws.Range("G5").Value = 0.8797879441049318
$ws.Range("H5").Value = 0.00852759115470823
$ws.Range("I5").Value = 0.8828324597582631
$ws.Range("J5").Value = 0.01524934530385911
$ws.Range("K5").Value = 0.8880002262187536
$ws.Range("L5").Value = 0.02033642124684555

# --- Row 6: RTREE ---------------------------------------------------------
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.9038574190302164
$ws.Range("D6").Value = 0.007997266537850396
$ws.Range("E6").Value = 0.8924296364431272
$ws.Range("F6").Value = 0.005975629623415342
$ws.Range("G6").Value = 0.8825875962916243
$ws.Range("H6").Value = 0.01109736525574546
$ws.Range("I6").Value = 0.8678807212006646
$ws.Range("J6").Value = 0.01291881660074068
$ws.Range("K6").Value = 0.8570851713607057
$ws.Range("L6").Value = 0.01574387731702154

# --- Row 7: XTREE ----------------------------------------------------------
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.9167100977198697
$ws.Range("D7").Value = 0.007272269145301402
$ws.Range("E7").Value = 0.9055808222755946
$ws.Range("F7").Value = 0.01047715717062251
$ws.Range("G7").Value = 0.8969379569642267
$ws.Range("H7").Value = 0.009588859623350775
$ws.Range("I7").Value = 0.8904053531534627
$ws.Range("J7").Value = 0.01821420338589705
$ws.Range("K7").Value = 0.889892546092071
$ws.Range("L7").Value = 0.01345068911210603

# --- Row 8: SVM (was row 9 before; the NB row that used to be row 8 is gone) -
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.9116654749609385
$ws.Range("D8").Value = 0.008471631091140162
$ws.Range("E8").Value = 0.9049071377064631
$ws.Range("F8").Value = 0.009946442248727021
$ws.Range("G8").Value = 0.9000881780216281
$ws.Range("H8").Value = 0.01280916016013744
$ws.Range("I8").Value = 0.8972188806782381
$ws.Range("J8").Value = 0.01513700033270456
$ws.Range("K8").Value = 0.8939424273272254
$ws.Range("L8").Value = 0.01173971300751821

# The table shrank by one row (8 rows total incl. header) - clear the old row 9
$ws.Range("A9:L9").Clear()
